$d = $word.ActiveDocument
$x = $d.WordOpenXML
Write-Host $x
